$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.996.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.44%  "

# Row 3
$ws.Range("D3").Value = "'3.535.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").Value = "'582.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

# Row 6
$ws.Range("D6").Value = "'179.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7
$ws.Range("D7").Value = "'3.528.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.96%  "

# Row 8
$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.44%  "

# Row 9
$ws.Range("E9").Value = "  +0.36%  "

# Row 10
$ws.Range("D10").Value = "'0.662"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.21%  "

# Row 11
$ws.Range("D11").Value = "'0.142"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.82%  "

# Row 12
$ws.Range("D12").Value = "'52.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.62%  "

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -15.25%  "

# Row 14
$ws.Range("D14").Value = "'9.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.29%  "

# Row 15
$ws.Range("D15").Value = "'4.119.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.12%  "

# Row 16
$ws.Range("D16").Value = "'3.542.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.45%  "

# Row 17
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("D18").Value = "'18.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.53%  "

# Row 19
$ws.Range("D19").Value = "'65.720.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.59%  "

# Row 20
$ws.Range("D20").Value = "'11.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.35%  "

# Row 21
$ws.Range("D21").Value = "'1.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.35%  "

# Row 22
$ws.Range("D22").Value = "'387.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.34%  "

# Row 23
$ws.Range("D23").Value = "'4.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.45%  "

# Row 24
$ws.Range("D24").Value = "'84.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.13%  "

# Row 25
$ws.Range("D25").Value = "'2.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.90%  "

# Row 26
$ws.Range("D26").Value = "'12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.72%  "

# Row 27
$ws.Range("D27").Value = "'5.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "

# Row 28
$ws.Range("D28").Value = "'10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.48%  "

# Row 29
$ws.Range("D29").Value = "'3.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.00%  "

# Row 30
$ws.Range("D30").Value = "'8.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.94%  "

# Row 31
$ws.Range("D31").Value = "'30.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.27%  "

# Row 32
$ws.Range("D32").Value = "'6.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.60%  "

# Row 33
$ws.Range("D33").Value = "'64.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "

# Row 34
$ws.Range("D34").Value = "'11.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "

# Row 35
$ws.Range("D35").Value = "'591.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.20%  "

# Row 36
$ws.Range("D36").Value = "'0.110"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.24%  "

# Row 37
$ws.Range("D37").Value = "'40.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.86%  "

# Row 38
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("E39").Value = "  -0.03%  "

# Row 40
$ws.Range("D40").Value = "'0.366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.48%  "

# Row 41
$ws.Range("D41").Value = "'0.0₃0725"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -17.96%  "

# Row 42
$ws.Range("D42").Value = "'0.127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.19%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.840.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.41%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.68%  "

# Row 45
$ws.Range("D45").Value = "'0.0403"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.76%  "

# Row 46
$ws.Range("D46").Value = "'2.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.07%  "

# Row 47
$ws.Range("D47").Value = "'0.129"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.31%  "

# Row 48
$ws.Range("D48").Value = "'3.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.74%  "

# Row 49
$ws.Range("D49").Value = "'2.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.64%  "

# Row 50
$ws.Range("D50").Value = "'134.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51
$ws.Range("D51").Value = "'8.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.26%  "
